# Auto-generated Excel COM-interop script to apply scheduled-runner data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for several
# rows across the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 38738.54
$ws.Range("I12").Value = 304.63635
$ws.Range("K12").Value = 304.63635
$ws.Range("M12").Value = -134.63635
# Row 40
$ws.Range("H40").Value = 1318.0358
$ws.Range("I40").Value = 1291.1666
$ws.Range("J40").Value = 1366.4
$ws.Range("K40").Value = 1291.1666
$ws.Range("L40").Value = 1366.4
$ws.Range("M40").Value = -1116.1666
$ws.Range("N40").Value = -1716.4
# Row 43
$ws.Range("H43").Value = 1496.0667
$ws.Range("I43").Value = 1941.5714
$ws.Range("J43").Value = 1106.25
$ws.Range("K43").Value = 1941.5714
$ws.Range("L43").Value = 1106.25
$ws.Range("M43").Value = -1872.5714
$ws.Range("N43").Value = -1244.25
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
# Row 106
$ws.Range("H106").Value = 1866.6666
$ws.Range("I106").Value = 1866.6666
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1866.6666
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1235.6666
$ws.Range("N106").ClearContents()
# Row 112
$ws.Range("H112").Value = 1167.091
$ws.Range("J112").Value = 1322
$ws.Range("L112").Value = 3966
$ws.Range("N112").Value = -6182
# Row 113
$ws.Range("H113").Value = 2188
$ws.Range("I113").Value = 2577.8
$ws.Range("J113").Value = 1863.1666
$ws.Range("K113").Value = 2577.8
$ws.Range("L113").Value = 1863.1666
$ws.Range("M113").Value = 676.1999999999998
$ws.Range("N113").Value = -8371.1666
# Row 132
$ws.Range("H132").Value = 4468704
$ws.Range("I132").Value = 4906541.5
$ws.Range("K132").Value = 14719624.5
$ws.Range("M132").Value = -14717094.5
# Row 135
$ws.Range("H135").Value = 1702.579
$ws.Range("I135").Value = 709.4167
$ws.Range("J135").Value = 3405.1428
$ws.Range("K135").Value = 6384.7503
$ws.Range("L135").Value = 30646.2852
$ws.Range("M135").Value = -3849.7503
$ws.Range("N135").Value = -35716.2852

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7879.517
$ws.Range("I32").Value = 7115.34
$ws.Range("J32").Value = 15979.8
$ws.Range("K32").Value = 7115.34
$ws.Range("L32").Value = 15979.8
$ws.Range("M32").Value = -6828.34
$ws.Range("N32").Value = -16553.8
# Row 102
$ws.Range("H102").Value = 144840
$ws.Range("I102").Value = 251972.5
$ws.Range("J102").Value = 1996.6666
$ws.Range("K102").Value = 251972.5
$ws.Range("L102").Value = 1996.6666
$ws.Range("M102").Value = -250350.5
$ws.Range("N102").Value = -5240.6666
# Row 122
$ws.Range("H122").Value = 1914.9
$ws.Range("I122").Value = 1749.1765
$ws.Range("K122").Value = 5247.529500000001
$ws.Range("M122").Value = -2797.529500000001
# Row 132
$ws.Range("H132").Value = 1749.3143
$ws.Range("I132").Value = 1749.3143
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5247.9429
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2717.9429
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 438
$ws.Range("I22").Value = 438
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 438
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -88
$ws.Range("N22").ClearContents()
# Row 62
$ws.Range("H62").Value = 2637.375
$ws.Range("J62").Value = 2650
$ws.Range("L62").Value = 2650
$ws.Range("N62").Value = -3898
# Row 65
$ws.Range("H65").Value = 2637.375
$ws.Range("J65").Value = 2650
$ws.Range("L65").Value = 13250
$ws.Range("N65").Value = -19490
# Row 105
$ws.Range("H105").Value = 1469.5385
$ws.Range("I105").Value = 1453
$ws.Range("J105").Value = 1496
$ws.Range("K105").Value = 1453
$ws.Range("L105").Value = 1496
$ws.Range("M105").Value = 294
$ws.Range("N105").Value = -4990

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1103.8372
$ws.Range("I5").Value = 973.35297
$ws.Range("J5").Value = 1189.1538
$ws.Range("K5").Value = 2920.05891
$ws.Range("L5").Value = 3567.4614
$ws.Range("M5").Value = -2808.05891
$ws.Range("N5").Value = -3791.4614
# Row 34
$ws.Range("H34").Value = 1379.6
$ws.Range("J34").Value = 1379.6
$ws.Range("L34").Value = 4138.799999999999
$ws.Range("N34").Value = -4306.799999999999
# Row 122
$ws.Range("H122").Value = 558.6667
$ws.Range("I122").Value = 539.8333
$ws.Range("J122").Value = 577.5
$ws.Range("K122").Value = 4858.4997
$ws.Range("L122").Value = 5197.5
$ws.Range("M122").Value = -2408.4997
$ws.Range("N122").Value = -10097.5
# Row 131
$ws.Range("H131").Value = 852.83
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 858.4141
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2575.2423
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12655.2423
# Row 132
$ws.Range("H132").Value = 1714.4286
# Row 135
$ws.Range("H135").Value = 1103.8372
$ws.Range("I135").Value = 973.35297
$ws.Range("J135").Value = 1189.1538
$ws.Range("K135").Value = 8760.176730000001
$ws.Range("L135").Value = 10702.3842
$ws.Range("M135").Value = -6225.176730000001
$ws.Range("N135").Value = -15772.3842

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 167.08333
$ws.Range("I2").Value = 86.71429000000001
$ws.Range("J2").Value = 279.6
$ws.Range("K2").Value = 86.71429000000001
$ws.Range("L2").Value = 279.6
$ws.Range("M2").Value = 26.28570999999999
$ws.Range("N2").Value = -505.6
# Row 70
$ws.Range("H70").Value = 44543.44
$ws.Range("I70").Value = 59911.11
$ws.Range("K70").Value = 59911.11
$ws.Range("M70").Value = -59641.11
# Row 73
$ws.Range("H73").Value = 44543.44
$ws.Range("I73").Value = 59911.11
$ws.Range("K73").Value = 59911.11
$ws.Range("M73").Value = -58975.11
# Row 102
$ws.Range("H102").Value = 2681
$ws.Range("I102").Value = 1949.5625
$ws.Range("J102").Value = 3851.3
$ws.Range("K102").Value = 1949.5625
$ws.Range("L102").Value = 3851.3
$ws.Range("M102").Value = -327.5625
$ws.Range("N102").Value = -7095.3

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3473.625
$ws.Range("I22").Value = 3697.25
$ws.Range("K22").Value = 3697.25
$ws.Range("M22").Value = -3402.25
# Row 27
$ws.Range("H27").Value = 3473.625
$ws.Range("I27").Value = 3697.25
$ws.Range("K27").Value = 3697.25
$ws.Range("M27").Value = -3590.25
# Row 122
$ws.Range("H122").Value = 2720
$ws.Range("I122").Value = 2720
$ws.Range("K122").Value = 8160
$ws.Range("M122").Value = -5710
# Row 123
$ws.Range("H123").Value = 23442.25
$ws.Range("J123").Value = 32494.5
$ws.Range("L123").Value = 32494.5
$ws.Range("N123").Value = -42294.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 111985.78
$ws.Range("I100").Value = 125546.5
$ws.Range("K100").Value = 251093
$ws.Range("M100").Value = -250552
# Row 126
$ws.Range("H126").Value = 1534.9166
$ws.Range("I126").Value = 1352
$ws.Range("K126").Value = 4056
$ws.Range("M126").Value = -1586
# Row 132
$ws.Range("H132").Value = 3032.3462
$ws.Range("I132").Value = 2994.2083
$ws.Range("K132").Value = 8982.624899999999
$ws.Range("M132").Value = -6452.624899999999
# Row 136
$ws.Range("H136").Value = 1691.303
$ws.Range("I136").Value = 633.4761999999999
$ws.Range("K136").Value = 1900.4286
$ws.Range("M136").Value = 649.5714000000003
